$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").Value = "NSE:HDFCSILVER"
$ws.Range("C2").Value = "NSE:360ONE"
$ws.Range("E2").Value = "NSE:ADANIGREEN"
$ws.Range("F2").Value = "NSE:HINDALCO"

# --- Row 3 ---
$ws.Range("B3").Value = "NSE:HINDALCO"
$ws.Range("C3").Value = "NSE:AXISTECETF"
$ws.Range("E3").Value = "NSE:ONGC"

# --- Row 4 (E4 cleared) ---
$ws.Range("C4").Value = "NSE:BALRAMCHIN"
$ws.Range("E4").ClearContents()

# --- Row 5 (E5 cleared) ---
$ws.Range("C5").Value = "NSE:BVCL"
$ws.Range("E5").ClearContents()

# --- Row 6 (E6 cleared) ---
$ws.Range("C6").Value = "NSE:DATAPATTNS"
$ws.Range("E6").ClearContents()

# --- Row 7 (E7 cleared) ---
$ws.Range("C7").Value = "NSE:DBL"
$ws.Range("E7").ClearContents()

# --- Row 8 (E8 cleared) ---
$ws.Range("C8").Value = "NSE:EICHERMOT"
$ws.Range("E8").ClearContents()

# --- Row 9 (E9 cleared) ---
$ws.Range("C9").Value = "NSE:ESCORTS"
$ws.Range("E9").ClearContents()

# --- Row 10 ---
$ws.Range("C10").Value = "NSE:EXCELINDUS"

# --- New rows 11-23: set A (index) with same format as A10, and C (ticker) ---
$newRows = @(
    @{ Row = 11; Idx = 9;  Ticker = "NSE:FINEORG" },
    @{ Row = 12; Idx = 10; Ticker = "NSE:FINOPB" },
    @{ Row = 13; Idx = 11; Ticker = "NSE:GULFPETRO" },
    @{ Row = 14; Idx = 12; Ticker = "NSE:IMPAL" },
    @{ Row = 15; Idx = 13; Ticker = "NSE:INDORAMA" },
    @{ Row = 16; Idx = 14; Ticker = "NSE:JKLAKSHMI" },
    @{ Row = 17; Idx = 15; Ticker = "NSE:M&M" },
    @{ Row = 18; Idx = 16; Ticker = "NSE:MAHESHWARI" },
    @{ Row = 19; Idx = 17; Ticker = "NSE:MONTECARLO" },
    @{ Row = 20; Idx = 18; Ticker = "NSE:MPHASIS" },
    @{ Row = 21; Idx = 19; Ticker = "NSE:PRINCEPIPE" },
    @{ Row = 22; Idx = 20; Ticker = "NSE:RATNAMANI" },
    @{ Row = 23; Idx = 21; Ticker = "NSE:RHIM" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Idx
    $ws.Range("C$r").Value = $item.Ticker

    # Copy the formatting from A10 (existing index-column style) onto the new A cell
    $ws.Range("A10").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
